# Fix "Nome da mãe" text to "Nome da mãe ou responsável" in the import
# template's header row (cell C1), per commit:
# "Correcao do texto Nome da mae ou responsavel no arquivo padrao xls para download"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the header text in C1
$ws.Range("C1").Value = "Nome da mãe ou responsável"

# Column C needs to widen to fit the new, longer text (AutoFit-style behavior).
# The workbook was saved with a best-fit width of ~25.29 characters for this
# column, so size it explicitly to match.
$ws.Columns("C").ColumnWidth = 24.45

# Move the active selection to C12, matching the saved workbook state
$ws.Range("C12").Select() | Out-Null
